# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (new row 147 and 148) for
# "Terminal Hortofrutícola Agro Chillán" - Ciruela, pushing the
# previously-existing rows 147-166 down to 149-168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 147, shifting
# everything from 147 down onward two rows lower (147->149 ... 166->168).
$ws.Rows("147:148").Insert()

# --- New row 147 ---
$ws.Range("A147").Value = 7
$ws.Range("B147").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C147").Value = "Ñuble"
$ws.Range("D147").Value2 = 45077
$ws.Range("E147").Value = 16
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100103
$ws.Range("H147").Value = "Frutos de hueso (carozo)"
$ws.Range("I147").Value = 100103002
$ws.Range("J147").Value = "Ciruela"
$ws.Range("K147").Value = "Angeleno"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 30
$ws.Range("N147").Value = 8000
$ws.Range("O147").Value = 8000
$ws.Range("P147").Value = 8000
$ws.Range("Q147").Value = "$/caja 18 kilos granel"
$ws.Range("R147").Value = "Región de O'Higgins"
$ws.Range("S147").Value = 444
$ws.Range("T147").Value = 18

# --- New row 148 ---
$ws.Range("A148").Value = 7
$ws.Range("B148").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C148").Value = "Ñuble"
$ws.Range("D148").Value2 = 45077
$ws.Range("E148").Value = 16
$ws.Range("F148").Value = "Fruta"
$ws.Range("G148").Value = 100103
$ws.Range("H148").Value = "Frutos de hueso (carozo)"
$ws.Range("I148").Value = 100103002
$ws.Range("J148").Value = "Ciruela"
$ws.Range("K148").Value = "Angeleno"
$ws.Range("L148").Value = "Segunda"
$ws.Range("M148").Value = 30
$ws.Range("N148").Value = 6000
$ws.Range("O148").Value = 6000
$ws.Range("P148").Value = 6000
$ws.Range("Q148").Value = "$/caja 18 kilos granel"
$ws.Range("R148").Value = "Región de O'Higgins"
$ws.Range("S148").Value = 333
$ws.Range("T148").Value = 18
